# Weekly fruit/vegetable price update: insert two new price records for
# "Ajo" (Chino variety, "Primera" quality) right above the existing block
# that starts at row 267, shifting the remaining rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 267-268 (existing rows 267.. shift to 269..)
$ws.Range("A267:A268").EntireRow.Insert()

# --- New row 267 ---
$ws.Cells.Item(267, 1).Value = 9
$ws.Cells.Item(267, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(267, 3).Value = "Metropolitana"
$ws.Cells.Item(267, 4).Value = 44858
$ws.Cells.Item(267, 5).Value = 13
$ws.Cells.Item(267, 6).Value = 100112003
$ws.Cells.Item(267, 7).Value = "Ajo"
$ws.Cells.Item(267, 8).Value = "Chino"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 400
$ws.Cells.Item(267, 11).Value = 13000
$ws.Cells.Item(267, 12).Value = 14000
$ws.Cells.Item(267, 13).Value = 13500
$ws.Cells.Item(267, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(267, 15).Value = "China"
$ws.Cells.Item(267, 16).Value = 1350
$ws.Cells.Item(267, 17).Value = 10
$ws.Cells.Item(267, 18).Value = "Hortaliza"

# --- New row 268 ---
$ws.Cells.Item(268, 1).Value = 9
$ws.Cells.Item(268, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(268, 3).Value = "Metropolitana"
$ws.Cells.Item(268, 4).Value = 44858
$ws.Cells.Item(268, 5).Value = 13
$ws.Cells.Item(268, 6).Value = 100112003
$ws.Cells.Item(268, 7).Value = "Ajo"
$ws.Cells.Item(268, 8).Value = "Chino"
$ws.Cells.Item(268, 9).Value = "Primera"
$ws.Cells.Item(268, 10).Value = 150
$ws.Cells.Item(268, 11).Value = 16000
$ws.Cells.Item(268, 12).Value = 16000
$ws.Cells.Item(268, 13).Value = 16000
$ws.Cells.Item(268, 14).Value = "`$/malla 10 kilos"
$ws.Cells.Item(268, 15).Value = "China"
$ws.Cells.Item(268, 16).Value = 1600
$ws.Cells.Item(268, 17).Value = 10
$ws.Cells.Item(268, 18).Value = "Hortaliza"
